# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The diff shows the serial date value in column C changing from
# 45725 (2025-03-09) to 45726 (2025-03-10) for every data row (2 through 43),
# i.e. every row was touched by the automatic update and its "changed" date
# bumped by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45726
}
